# Update the Training Dashboard sheet with the new progress date (04-Nov-2025).
# For rows 3 through 20:
#   - Column H ("PERIOD TO EXPIRE") decreases by 1 day
#   - Column I ("LAST UPDATE") changes from 03-Nov-2025 to 04-Nov-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 20; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H ("PERIOD TO EXPIRE")
    $hCell.Value = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I ("LAST UPDATE")
    # Force the cell to remain plain text (rather than being auto-converted
    # to a date serial value) when assigning the new date string.
    $iCell.NumberFormat = "@"
    $iCell.Value = "04-Nov-2025"

    # Restore the original (General) number formatting on column I by
    # copying the format from its row-mate in column J, which was never
    # touched, so the saved file keeps referencing the same style as
    # before the edit (handles the differently-styled rows 18-19 too).
    $formatSource = $ws.Cells.Item($row, 10)
    $formatSource.Copy()
    $iCell.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false
